# Insert two new weekly records at the top of the Rabanito price history
# (rows 67-68 of Sheet1), pushing the existing rows 67-171 down to 69-173.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 67..171 down by two rows -> 69..173.
$ws.Rows("67:68").Insert()

# New record for (new) row 67.
$ws.Range("A67").Value = 9
$ws.Range("B67").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C67").Value = "Metropolitana"
$ws.Range("D67").Value = 44495
$ws.Range("E67").Value = 13
$ws.Range("F67").Value = 300000001
$ws.Range("G67").Value = "Rabanito"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 15000
$ws.Range("K67").Value = 3500
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = 3733
$ws.Range("N67").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O67").Value = "Provincia de Chacabuco"
$ws.Range("P67").Value = 37
$ws.Range("Q67").Value = 100
$ws.Range("R67").Value = "Hortaliza"

# New record for (new) row 68.
$ws.Range("A68").Value = 9
$ws.Range("B68").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C68").Value = "Metropolitana"
$ws.Range("D68").Value = 44495
$ws.Range("E68").Value = 13
$ws.Range("F68").Value = 300000001
$ws.Range("G68").Value = "Rabanito"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Segunda"
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = 3000
$ws.Range("N68").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O68").Value = "Provincia de Chacabuco"
$ws.Range("P68").Value = 30
$ws.Range("Q68").Value = 100
$ws.Range("R68").Value = "Hortaliza"
